$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D2:D51 range to text format before writing values so that
# numeric-looking price strings (e.g. "587.29") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.301.06'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '3.495.99'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '587.29'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = '134.20'
$ws.Range('E6').Value = '  +1.58%  '
$ws.Range('D7').Value = '3.494.51'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').Value = '4.087.99'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D16').Value = '3.494.36'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').Value = '64.336.18'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '25.27'
$ws.Range('E18').Value = '  -8.92%  '
$ws.Range('D19').Value = '9.88'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').Value = '5.74'
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('D21').Value = '13.61'
$ws.Range('E21').Value = '  -6.30%  '
$ws.Range('D22').Value = '387.64'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('B23').Value = 'WrappedeETH'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D23').Value = '3.634.04'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').Value = '0.565'
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '7.42'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = '8.29'
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').Value = '3.515.45'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +2.94%  '
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').Value = '5.26'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').Value = '6.88'
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('D40').Value = '1.54'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').Value = '161.90'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D45').Value = '25.42'
$ws.Range('E45').Value = '  -6.82%  '
$ws.Range('D46').Value = '41.94'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').Value = '1.19'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').Value = '2.475.65'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('E51').Value = '  -2.01%  '

# Restore the original (default) cell style now that the text values are set,
# so no residual number-format styling is left behind.
$priceRange.Style = "Normal"
